$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:H1, matching style of existing header (e.g. E1) ---
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Boolean outlier flag data, rows 2-14 for columns F (KNN), G (SVM), H (RF) ---
$flagsF = @($false, $false, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
$flagsG = @($false, $false, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true)
$flagsH = @($false, $false, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true)

for ($i = 0; $i -lt 13; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $flagsF[$i]
    $ws.Cells.Item($row, 7).Value = $flagsG[$i]
    $ws.Cells.Item($row, 8).Value = $flagsH[$i]
}
